# Auto-generated edit script applying the Halicarnassus_Profits diff
# Updates literal (non-formula) market-data values across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()

$ws.Cells.Item(88, 8).Value = 1036.4445
$ws.Cells.Item(88, 9).Value = 1174.75
$ws.Cells.Item(88, 10).Value = 925.8
$ws.Cells.Item(88, 11).Value = 1174.75
$ws.Cells.Item(88, 12).Value = 925.8
$ws.Cells.Item(88, 13).Value = -768.75
$ws.Cells.Item(88, 14).Value = -1737.8

$ws.Cells.Item(91, 8).Value = 1036.4445
$ws.Cells.Item(91, 9).Value = 1174.75
$ws.Cells.Item(91, 10).Value = 925.8
$ws.Cells.Item(91, 11).Value = 1174.75
$ws.Cells.Item(91, 12).Value = 925.8
$ws.Cells.Item(91, 13).Value = 229.25
$ws.Cells.Item(91, 14).Value = -3733.8

$ws.Cells.Item(99, 8).Value = 3727.4666
$ws.Cells.Item(99, 9).Value = 1859.3334
$ws.Cells.Item(99, 10).Value = 4194.5
$ws.Cells.Item(99, 11).Value = 5578.0002
$ws.Cells.Item(99, 12).Value = 12583.5
$ws.Cells.Item(99, 13).Value = -4080.0002
$ws.Cells.Item(99, 14).Value = -15579.5

$ws.Cells.Item(125, 8).Value = 5340
$ws.Cells.Item(125, 9).Value = 4500
$ws.Cells.Item(125, 10).Value = 5508
$ws.Cells.Item(125, 11).Value = 40500
$ws.Cells.Item(125, 12).Value = 49572
$ws.Cells.Item(125, 13).Value = -38040
$ws.Cells.Item(125, 14).Value = -54492

$ws.Cells.Item(137, 8).Value = 1596.7142
$ws.Cells.Item(137, 9).Value = 1135.4
$ws.Cells.Item(137, 10).Value = 2750
$ws.Cells.Item(137, 11).Value = 3406.2
$ws.Cells.Item(137, 12).Value = 8250
$ws.Cells.Item(137, 13).Value = -856.2000000000003
$ws.Cells.Item(137, 14).Value = -13350

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2475.1
$ws.Cells.Item(45, 10).Value = 3593
$ws.Cells.Item(45, 12).Value = 3593
$ws.Cells.Item(45, 14).Value = -4347

$ws.Cells.Item(61, 8).Value = 6848.9287
$ws.Cells.Item(61, 9).Value = 5697.857
$ws.Cells.Item(61, 11).Value = 5697.857
$ws.Cells.Item(61, 13).Value = -5485.857

$ws.Cells.Item(74, 8).Value = 3412.125
$ws.Cells.Item(74, 9).Value = 3628.4285
$ws.Cells.Item(74, 10).Value = 1898
$ws.Cells.Item(74, 11).Value = 3628.4285
$ws.Cells.Item(74, 12).Value = 1898
$ws.Cells.Item(74, 13).Value = -2754.4285
$ws.Cells.Item(74, 14).Value = -3646

$ws.Cells.Item(77, 8).Value = 3412.125
$ws.Cells.Item(77, 9).Value = 3628.4285
$ws.Cells.Item(77, 10).Value = 1898
$ws.Cells.Item(77, 11).Value = 18142.1425
$ws.Cells.Item(77, 12).Value = 9490
$ws.Cells.Item(77, 13).Value = -13774.1425
$ws.Cells.Item(77, 14).Value = -18226

$ws.Cells.Item(132, 8).Value = 2689.423
$ws.Cells.Item(132, 10).Value = 4540.5
$ws.Cells.Item(132, 12).Value = 13621.5
$ws.Cells.Item(132, 14).Value = -18681.5

$ws.Cells.Item(136, 8).Value = 6848.9287
$ws.Cells.Item(136, 9).Value = 5697.857
$ws.Cells.Item(136, 11).Value = 17093.571
$ws.Cells.Item(136, 13).Value = -14543.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4520.778
$ws.Cells.Item(107, 9).Value = 3644.2307
$ws.Cells.Item(107, 11).Value = 3644.2307
$ws.Cells.Item(107, 13).Value = -1724.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5795.9033
$ws.Cells.Item(31, 9).Value = 2677.2307
$ws.Cells.Item(31, 10).Value = 8048.278
$ws.Cells.Item(31, 11).Value = 2677.2307
$ws.Cells.Item(31, 12).Value = 8048.278
$ws.Cells.Item(31, 13).Value = -2382.2307
$ws.Cells.Item(31, 14).Value = -8638.278

$ws.Cells.Item(34, 8).Value = 5795.9033
$ws.Cells.Item(34, 9).Value = 2677.2307
$ws.Cells.Item(34, 10).Value = 8048.278
$ws.Cells.Item(34, 11).Value = 2677.2307
$ws.Cells.Item(34, 12).Value = 8048.278
$ws.Cells.Item(34, 13).Value = -2475.2307
$ws.Cells.Item(34, 14).Value = -8452.278

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 411.46667
$ws.Cells.Item(34, 9).Value = 89.8
$ws.Cells.Item(34, 10).Value = 572.3
$ws.Cells.Item(34, 11).Value = 269.4
$ws.Cells.Item(34, 12).Value = 1716.9
$ws.Cells.Item(34, 13).Value = -185.4
$ws.Cells.Item(34, 14).Value = -1884.9

$ws.Cells.Item(112, 8).Value = 710.5
$ws.Cells.Item(112, 9).Value = 713
$ws.Cells.Item(112, 10).Value = 708
$ws.Cells.Item(112, 11).Value = 2139
$ws.Cells.Item(112, 12).Value = 2124
$ws.Cells.Item(112, 13).Value = -1031
$ws.Cells.Item(112, 14).Value = -4340

$ws.Cells.Item(128, 8).Value = 605994.4
$ws.Cells.Item(128, 9).Value = 605994.4
$ws.Cells.Item(128, 11).Value = 1817983.2
$ws.Cells.Item(128, 13).Value = -1813003.2

$ws.Cells.Item(129, 8).Value = 1997.6666
$ws.Cells.Item(129, 10).Value = 2000
$ws.Cells.Item(129, 12).Value = 6000
$ws.Cells.Item(129, 14).Value = -16000

$ws.Cells.Item(132, 8).Value = 2004.0416
$ws.Cells.Item(132, 9).Value = 1566.6666
$ws.Cells.Item(132, 10).Value = 2066.524
$ws.Cells.Item(132, 11).Value = 14099.9994
$ws.Cells.Item(132, 12).Value = 18598.716
$ws.Cells.Item(132, 13).Value = -11569.9994
$ws.Cells.Item(132, 14).Value = -23658.716

$ws.Cells.Item(136, 8).Value = 10666.667
$ws.Cells.Item(136, 9).Value = 9000
$ws.Cells.Item(136, 11).Value = 27000
$ws.Cells.Item(136, 13).Value = -21900

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1664.2632
$ws.Cells.Item(102, 9).Value = 1271.8823
$ws.Cells.Item(102, 11).Value = 1271.8823
$ws.Cells.Item(102, 13).Value = 350.1177

$ws.Cells.Item(126, 8).Value = 2943
$ws.Cells.Item(126, 9).Value = 2943
$ws.Cells.Item(126, 11).Value = 8829
$ws.Cells.Item(126, 13).Value = -6359

$ws.Cells.Item(136, 8).Value = 49999.5
$ws.Cells.Item(136, 10).Value = 49999.5
$ws.Cells.Item(136, 12).Value = 149998.5
$ws.Cells.Item(136, 14).Value = -155098.5

$ws.Cells.Item(141, 8).Value = 23997.5
$ws.Cells.Item(141, 10).Value = 23997.5
$ws.Cells.Item(141, 12).Value = 23997.5
$ws.Cells.Item(141, 14).Value = -34357.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 7777
$ws.Cells.Item(16, 10).Value = 7777
$ws.Cells.Item(16, 12).Value = 7777
$ws.Cells.Item(16, 14).Value = -8117

$ws.Cells.Item(61, 8).Value = 3610.25
$ws.Cells.Item(61, 9).Value = 976.4
$ws.Cells.Item(61, 11).Value = 976.4
$ws.Cells.Item(61, 13).Value = -774.4

$ws.Cells.Item(113, 8).Value = 3610.25
$ws.Cells.Item(113, 9).Value = 976.4
$ws.Cells.Item(113, 11).Value = 976.4
$ws.Cells.Item(113, 13).Value = 1193.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3888.5417
$ws.Cells.Item(136, 9).Value = 2577.6667
$ws.Cells.Item(136, 10).Value = 5199.4165
$ws.Cells.Item(136, 11).Value = 7733.000100000001
$ws.Cells.Item(136, 12).Value = 15598.2495
$ws.Cells.Item(136, 13).Value = -5183.000100000001
$ws.Cells.Item(136, 14).Value = -20698.2495
